# Daily attendance processing - 2026-01-19 11:13:00
#
# The automated attendance sync re-wrote the "Recorded By" audit trail for
# rows that were re-processed by the backup/system account: the trailing
# "System"/"system"/"backup@backdoor.com" entry now comes first in the
# comma-separated list (e.g. "dnasr281@gmail.com, System" becomes
# "System, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

# Exact "Recorded By" strings that were reprocessed, mapped old -> new.
$map = @{
    "system, System, backup@backdoor.com" = "backup@backdoor.com, system, System";
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com";
    "System, backup@backdoor.com"         = "backup@backdoor.com, System";
}

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
